$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The underlying records for rows 15-18 get re-sorted (same 4 observations,
# different row order). Only columns A,B,D,E,F,G,H,I,J,P,Q,R (+ the L / AC
# cells that only exist on some rows) actually change value; every other
# column (C,N,S,T,U,V,W,Y,Z,AA,AB,AD,AE,AF,AG,AT,AW,AX,AY) is identical
# across the four rows and is left untouched.
# ---------------------------------------------------------------------------

# ---- Row 15 (becomes the former row16 record) -----------------------------
$ws.Range("A15").Value = 111837758
$ws.Range("B15").Value = 90187
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 2014
$ws.Range("F15").Value = "Koralltaggsvamp"
$ws.Range("G15").Value = "Hericium coralloides"
$ws.Range("H15").Value = "(Scop.:Fr.) Pers."
$ws.Range("I15").Value = "'6"
$ws.Range("I15").ClearFormats()
$ws.Range("J15").Value = "fruktkroppar"
$ws.Range("P15").Value = "Brotorp, hyggeskant, Sm"
$ws.Range("Q15").Value = 575673.5681218
$ws.Range("R15").Value = 6404513.458820416
$ws.Range("L15").ClearContents()
$ws.Range("AC15").Value = "På asplåga."

# ---- Row 16 (becomes the former row18 record) ------------------------------
$ws.Range("A16").Value = 111837705
$ws.Range("B16").Value = 90662
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 4363
$ws.Range("F16").Value = "Zontaggsvamp"
$ws.Range("G16").Value = "Hydnellum concrescens"
$ws.Range("H16").Value = "(Pers.) Banker"
$ws.Range("I16").Value = "'10"
$ws.Range("I16").ClearFormats()
$ws.Range("J16").Value = "fruktkroppar"
$ws.Range("P16").Value = "Brotorp, Långsjön, Sm"
$ws.Range("Q16").Value = 575795.3141537429
$ws.Range("R16").Value = 6404518.948622406
$ws.Range("AC16").ClearContents()

# ---- Row 17 (becomes the former row15 record) ------------------------------
$ws.Range("A17").Value = 111837675
$ws.Range("B17").Value = 103288
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 221144
$ws.Range("F17").Value = "Grönpyrola"
$ws.Range("G17").Value = "Pyrola chlorantha"
$ws.Range("H17").Value = "Sw."
$ws.Range("I17").Value = "'10"
$ws.Range("I17").ClearFormats()
$ws.Range("J17").Value = "plantor/tuvor"
$ws.Range("P17").Value = "Brotorp, Långsjön, Sm"
$ws.Range("Q17").Value = 575781.9606960951
$ws.Range("R17").Value = 6404546.96767282
$ws.Range("L17").Value = "'"
$ws.Range("L17").ClearFormats()

# ---- Row 18 (becomes the former row17 record) ------------------------------
$ws.Range("A18").Value = 111837741
$ws.Range("B18").Value = 90658
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 4361
$ws.Range("F18").Value = "Orange taggsvamp"
$ws.Range("G18").Value = "Hydnellum aurantiacum"
$ws.Range("H18").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I18").Value = "'15"
$ws.Range("I18").ClearFormats()
$ws.Range("J18").Value = "fruktkroppar"
$ws.Range("P18").Value = "Brotorp, hyggeskant, Sm"
$ws.Range("Q18").Value = 575653.9215098171
$ws.Range("R18").Value = 6404506.688862759
